$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 5)
$ws.Range("D2").Value = 0.9999589069136334
$ws.Range("E2").Value = 0.9999589069136334

# Row 3 (MDD 27)
$ws.Range("C3").Value = $false
$ws.Range("D3").Value = 0.9999989610682003
$ws.Range("E3").Value = 0.9999989610682003

# Row 4 (MDD 47)
$ws.Range("D4").Value = 0.9999999999999667
$ws.Range("E4").Value = [double]"3.33066907387547E-14"

# Row 5 (MDD 13)
$ws.Range("C5").Value = $false
$ws.Range("D5").Value = 0.0003343057064040079
$ws.Range("E5").Value = 0.999665694293596

# Row 6 (MDD 25)
$ws.Range("C6").Value = $true
$ws.Range("D6").Value = 0.6791890421223865
$ws.Range("E6").Value = 0.3208109578776135

# Row 8 (MDD 5)
$ws.Range("D8").Value = 0.003856858119883094
$ws.Range("E8").Value = 0.9961431418801169
$ws.Range("F8").Value = 5.403600215911865
$ws.Range("G8").Value = 0.4285714285714285
